# Clean_Actions_Report_Modeling.xlsx — "finished ls model" update.
#
# Adds the two new rows produced once the LR (Logistic Regression) and RF
# (Random Forest) final modeling data frames finished their feature-pruning
# pass, and updates the "CompleteDF" rows' ResultShape now that gsubind/low
# correlation columns were dropped (348, 163) -> (348, 161).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the ResultShape for the two existing CompleteDF actions ---
$ws.Range("E2").Value = "(348, 161)"
$ws.Range("E3").Value = "(348, 161)"

# --- New row: LR_Final_ModelingDF feature pruning ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "LR_Final_ModelingDF"
$ws.Range("C4").Value = "Drop features w/ low impDecrease from Round 2"
$ws.Range("D4").Value = "['at', 'lse', 'pi_std', 'rest_count_of_diffs', 'sstk_std', 'xsga']"
$ws.Range("E4").Value = "(347, 23)"

# --- New row: RF_Final_ModelingDF feature pruning ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "RF_Final_ModelingDF"
$ws.Range("C5").Value = "Drop features w/ low impDecrease from Round 2"
$ws.Range("D5").Value = "['rat_spcsrc']"
$ws.Range("E5").Value = "(347, 28)"

# Match the existing index-column style (bold/bordered/centered) on the new rows.
$ws.Range("A2").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Widen the columns to fit the new (longer) content, matching the author's
# manual column sizing.
$ws.Columns.Item(2).ColumnWidth = 19
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 20
$ws.Columns.Item(5).ColumnWidth = 10.8333333333

# Leave the selection where the author left off.
$ws.Range("E5").Select()
